$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    'E1' = 'Saida (0,2*avg - 10) * 3'
    'B2' = '50,3'
    'C2' = '44,1'
    'D2' = '47,20'
    'E2' = '-1,68'
    'B3' = '54'
    'C3' = '54'
    'D3' = '54,00'
    'E3' = '2,40'
    'B4' = '47,1'
    'C4' = '48,6'
    'D4' = '47,85'
    'E4' = '-1,29'
    'B5' = '48'
    'C5' = '46,8'
    'D5' = '47,40'
    'E5' = '-1,56'
    'B6' = '50,1'
    'C6' = '47,4'
    'D6' = '48,75'
    'E6' = '-0,75'
    'B7' = '45,1'
    'C7' = '48'
    'D7' = '46,55'
    'E7' = '-2,07'
    'B8' = '51,8'
    'C8' = '50,6'
    'D8' = '51,20'
    'E8' = '0,72'
    'B9' = '48,8'
    'C9' = '54,1'
    'D9' = '51,45'
    'E9' = '0,87'
    'B10' = '55,9'
    'C10' = '53'
    'D10' = '54,45'
    'E10' = '2,67'
    'B11' = '49,4'
    'C11' = '49,5'
    'D11' = '49,45'
    'E11' = '-0,33'
    'B12' = '50,7'
    'C12' = '49,2'
    'D12' = '49,95'
    'E12' = '-0,03'
    'B13' = '43,8'
    'C13' = '46,1'
    'D13' = '44,95'
    'E13' = '-3,03'
    'B14' = '49,4'
    'C14' = '54,3'
    'D14' = '51,85'
    'E14' = '1,11'
    'B15' = '58,4'
    'C15' = '58,7'
    'D15' = '58,55'
    'E15' = '5,13'
    'B16' = '52,6'
    'C16' = '30,3'
    'D16' = '41,45'
    'E16' = '-5,13'
    'B17' = '49,9'
    'C17' = '58,1'
    'D17' = '54,00'
    'E17' = '2,40'
    'B18' = '53,6'
    'C18' = '55,9'
    'D18' = '54,75'
    'E18' = '2,85'
    'B19' = '42,8'
    'C19' = '45,3'
    'D19' = '44,05'
    'E19' = '-3,57'
    'B20' = '50,9'
    'C20' = '54,4'
    'D20' = '52,65'
    'E20' = '1,59'
    'B21' = '43,1'
    'C21' = '55,4'
    'D21' = '49,25'
    'E21' = '-0,45'
    'B22' = '54,9'
    'C22' = '62,6'
    'D22' = '58,75'
    'E22' = '5,25'
    'B23' = '49,8'
    'C23' = '46,9'
    'D23' = '48,35'
    'E23' = '-0,99'
    'B24' = '54,8'
    'C24' = '47,2'
    'D24' = '51,00'
    'E24' = '0,60'
    'B25' = '44,6'
    'C25' = '51,2'
    'D25' = '47,90'
    'E25' = '-1,26'
    'B26' = '42,2'
    'C26' = '42,3'
    'D26' = '42,25'
    'E26' = '-4,65'
    'B27' = '58'
    'C27' = '60,7'
    'D27' = '59,35'
    'E27' = '5,61'
    'B28' = '37,7'
    'C28' = '36,5'
    'D28' = '37,10'
    'E28' = '-7,74'
    'B29' = '48,3'
    'C29' = '48,6'
    'D29' = '48,45'
    'E29' = '-0,93'
    'B30' = '40,5'
    'C30' = '44,6'
    'D30' = '42,55'
    'E30' = '-4,47'
    'B31' = '42,6'
    'C31' = '50,2'
    'D31' = '46,40'
    'E31' = '-2,16'
    'B32' = '53,1'
    'C32' = '52,7'
    'D32' = '52,90'
    'E32' = '1,74'
    'B33' = '47,4'
    'C33' = '35,6'
    'D33' = '41,50'
    'E33' = '-5,10'
    'B34' = '53,4'
    'C34' = '50,3'
    'D34' = '51,85'
    'E34' = '1,11'
    'B35' = '47,7'
    'C35' = '55,2'
    'D35' = '51,45'
    'E35' = '0,87'
    'B36' = '53,1'
    'C36' = '41,4'
    'D36' = '47,25'
    'E36' = '-1,65'
    'B37' = '44,8'
    'C37' = '51,7'
    'D37' = '48,25'
    'E37' = '-1,05'
    'B38' = '48,8'
    'C38' = '49,1'
    'D38' = '48,95'
    'E38' = '-0,63'
    'B39' = '58,9'
    'C39' = '55,5'
    'D39' = '57,20'
    'E39' = '4,32'
    'B40' = '46,6'
    'C40' = '69,6'
    'D40' = '58,10'
    'E40' = '4,86'
    'B41' = '45,9'
    'C41' = '44,4'
    'D41' = '45,15'
    'E41' = '-2,91'
    'B42' = '54,7'
    'C42' = '54,5'
    'D42' = '54,60'
    'E42' = '2,76'
    'B43' = '55,3'
    'C43' = '49,7'
    'D43' = '52,50'
    'E43' = '1,50'
    'B44' = '53,1'
    'C44' = '53,4'
    'D44' = '53,25'
    'E44' = '1,95'
    'B45' = '47,9'
    'C45' = '30,5'
    'D45' = '39,20'
    'E45' = '-6,48'
    'B46' = '51,1'
    'C46' = '51,1'
    'D46' = '51,10'
    'E46' = '0,66'
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
    $cell.Style = "Normal"
}